$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 451.8889
$ws.Range("I28").Value = 509.06668
$ws.Range("J28").Value = 166
$ws.Range("K28").Value = 509.06668
$ws.Range("L28").Value = 166
$ws.Range("M28").Value = -24.06668000000002
$ws.Range("N28").Value = -1136

$ws.Range("H51").Value = 1937
$ws.Range("I51").Value = 1699
$ws.Range("J51").Value = 1976.6666
$ws.Range("K51").Value = 1699
$ws.Range("L51").Value = 1976.6666
$ws.Range("M51").Value = -1215
$ws.Range("N51").Value = -2944.6666

$ws.Range("H87").Value = 16831.578
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 16831.578
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 16831.578
$ws.Range("N87").Value = -19327.578

$ws.Range("H90").Value = 16831.578
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 16831.578
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 50494.734
$ws.Range("N90").Value = -62974.734

$ws.Range("H111").Value = 1165.4286
$ws.Range("I111").Value = 349
$ws.Range("J111").Value = 1492
$ws.Range("K111").Value = 1047
$ws.Range("L111").Value = 4476
$ws.Range("M111").Value = 2020
$ws.Range("N111").Value = -10610

$ws.Range("H115").Value = 2087.5
$ws.Range("I115").Value = 1000
$ws.Range("J115").Value = 2450
$ws.Range("K115").Value = 3000
$ws.Range("L115").Value = 7350
$ws.Range("M115").Value = -1433
$ws.Range("N115").Value = -10484

$ws.Range("H135").Value = 931.1094000000001
$ws.Range("I135").Value = 348.18
$ws.Range("J135").Value = 3013
$ws.Range("K135").Value = 3133.62
$ws.Range("L135").Value = 27117
$ws.Range("M135").Value = -598.6199999999999
$ws.Range("N135").Value = -32187

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()

$ws.Range("H13").Value = 60000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 60000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 60000
$ws.Range("N13").Value = -60288
$ws.Range("M13").ClearContents()

$ws.Range("H52").Value = 46945
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 46945
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 46945
$ws.Range("N52").Value = -47581

$ws.Range("H132").Value = 1440.3816
$ws.Range("I132").Value = 1338.017
$ws.Range("J132").Value = 1795.6471
$ws.Range("K132").Value = 4014.051
$ws.Range("L132").Value = 5386.9413
$ws.Range("M132").Value = -1484.051
$ws.Range("N132").Value = -10446.9413

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 38528
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 38528
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 38528
$ws.Range("N27").Value = -38912

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H55").Value = 50000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 50000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 50000
$ws.Range("N55").Value = -50546

$ws.Range("H86").Value = 3968.8696
$ws.Range("I86").Value = 3655.625
$ws.Range("J86").Value = 4684.857
$ws.Range("K86").Value = 3655.625
$ws.Range("L86").Value = 4684.857
$ws.Range("M86").Value = -2532.625
$ws.Range("N86").Value = -6930.857

$ws.Range("H89").Value = 3968.8696
$ws.Range("I89").Value = 3655.625
$ws.Range("J89").Value = 4684.857
$ws.Range("K89").Value = 18278.125
$ws.Range("L89").Value = 23424.285
$ws.Range("M89").Value = -12662.125
$ws.Range("N89").Value = -34656.285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 703.6667
$ws.Range("I16").Value = 705.5
$ws.Range("J16").Value = 700
$ws.Range("K16").Value = 705.5
$ws.Range("L16").Value = 700
$ws.Range("M16").Value = -418.5
$ws.Range("N16").Value = -1274

$ws.Range("H31").Value = 22708.268
$ws.Range("I31").Value = 24777.795
$ws.Range("J31").Value = 17017.062
$ws.Range("K31").Value = 24777.795
$ws.Range("L31").Value = 17017.062
$ws.Range("M31").Value = -24482.795
$ws.Range("N31").Value = -17607.062

$ws.Range("H34").Value = 22708.268
$ws.Range("I34").Value = 24777.795
$ws.Range("J34").Value = 17017.062
$ws.Range("K34").Value = 24777.795
$ws.Range("L34").Value = 17017.062
$ws.Range("M34").Value = -24575.795
$ws.Range("N34").Value = -17421.062

$ws.Range("H113").Value = 703.6667
$ws.Range("I113").Value = 705.5
$ws.Range("J113").Value = 700
$ws.Range("K113").Value = 705.5
$ws.Range("L113").Value = 700
$ws.Range("M113").Value = 1464.5
$ws.Range("N113").Value = -5040

$ws.Range("H132").Value = 894.74286
$ws.Range("I132").Value = 731.9655
$ws.Range("J132").Value = 1681.5
$ws.Range("K132").Value = 2195.8965
$ws.Range("L132").Value = 5044.5
$ws.Range("M132").Value = 334.1035000000002
$ws.Range("N132").Value = -10104.5

$ws.Range("H134").Value = 1688.909
$ws.Range("I134").Value = 1423.8182
$ws.Range("J134").Value = 2484.182
$ws.Range("K134").Value = 4271.4546
$ws.Range("L134").Value = 7452.545999999999
$ws.Range("M134").Value = -1736.4546
$ws.Range("N134").Value = -12522.546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1815
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 1815
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 5445
$ws.Range("N75").Value = -7441
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 1815
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 1815
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 16335
$ws.Range("N78").Value = -26319
$ws.Range("M78").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1677.1538
$ws.Range("I113").Value = 1881
$ws.Range("J113").Value = 997.6667
$ws.Range("K113").Value = 1881
$ws.Range("L113").Value = 997.6667
$ws.Range("M113").Value = 289
$ws.Range("N113").Value = -5337.6667

$ws.Range("H132").Value = 1496.5
$ws.Range("I132").Value = 1400.9286
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 4202.7858
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -1672.7858
$ws.Range("N132").Value = -12560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 37502
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 37502
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 37502
$ws.Range("N12").Value = -37842

$ws.Range("H61").Value = 1982.8
$ws.Range("I61").Value = 1780.25
$ws.Range("J61").Value = 2214.2856
$ws.Range("K61").Value = 1780.25
$ws.Range("L61").Value = 2214.2856
$ws.Range("M61").Value = -1578.25
$ws.Range("N61").Value = -2618.2856

$ws.Range("H113").Value = 1982.8
$ws.Range("I113").Value = 1780.25
$ws.Range("J113").Value = 2214.2856
$ws.Range("K113").Value = 1780.25
$ws.Range("L113").Value = 2214.2856
$ws.Range("M113").Value = 389.75
$ws.Range("N113").Value = -6554.2856

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 32002.5
$ws.Range("I10").Value = 4005
$ws.Range("J10").Value = 60000
$ws.Range("K10").Value = 4005
$ws.Range("L10").Value = 60000
$ws.Range("M10").Value = -3836
$ws.Range("N10").Value = -60338

$ws.Range("H13").Value = 19900
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 19900
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 19900
$ws.Range("N13").Value = -20180

$ws.Range("H122").Value = 4004581
$ws.Range("I122").Value = 2030785.5
$ws.Range("J122").Value = 10419416
$ws.Range("K122").Value = 6092356.5
$ws.Range("L122").Value = 31258248
$ws.Range("M122").Value = -6089906.5
$ws.Range("N122").Value = -31263148

$ws.Range("H126").Value = 3614977.8
$ws.Range("I126").Value = 3607468.5
$ws.Range("J126").Value = 3624990
$ws.Range("K126").Value = 10822405.5
$ws.Range("L126").Value = 10874970
$ws.Range("M126").Value = -10819935.5
$ws.Range("N126").Value = -10879910

$ws.Range("H132").Value = 924.6177
$ws.Range("I132").Value = 798.48
$ws.Range("J132").Value = 1275
$ws.Range("K132").Value = 2395.44
$ws.Range("L132").Value = 3825
$ws.Range("M132").Value = 134.5599999999999
$ws.Range("N132").Value = -8885
